$wb = $excel.ActiveWorkbook

# --- Fix typo "cep290_unkown" -> "cep290_unknown" on the "genotype" sheet ---
# (every cell that previously held the misspelled string gets corrected)
$ws = $wb.Worksheets.Item("genotype")

$cellsToFix = @("I4", "J5", "K5", "F6", "J6", "K6", "L6", "M6", "D9", "E9", "I9")
foreach ($addr in $cellsToFix) {
    $ws.Range($addr).Value = "cep290_unknown"
}

# Make "genotype" the active/selected sheet, with D34 as the active cell
$ws.Activate()
$ws.Range("D34").Select()
